# Fruta / hortaliza, semanal
# Insert two new weekly price rows above the old row 25, pushing the
# existing data (old rows 25-105) down to rows 27-107.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 25 (each Insert() shifts everything below
# it down by one row, so calling it twice at the same index opens a
# 2-row gap at rows 25-26).
$ws.Rows.Item(25).Insert()
$ws.Rows.Item(25).Insert()

# New row 25 - Poroto verde, Sin especificar, Primera
$ws.Cells.Item(25, 1).Value = 1
$ws.Cells.Item(25, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(25, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(25, 4).Value = 45250
$ws.Cells.Item(25, 5).Value = 15
$ws.Cells.Item(25, 6).Value = 100112031
$ws.Cells.Item(25, 7).Value = "Poroto verde"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 2500
$ws.Cells.Item(25, 11).Value = 1800
$ws.Cells.Item(25, 12).Value = 1800
$ws.Cells.Item(25, 13).Value = 1800
$ws.Cells.Item(25, 14).Value = "`$/kilo"
$ws.Cells.Item(25, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(25, 16).Value = 1800
$ws.Cells.Item(25, 17).Value = 1
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# New row 26 - Poroto verde, Sin especificar, Segunda
$ws.Cells.Item(26, 1).Value = 1
$ws.Cells.Item(26, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(26, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(26, 4).Value = 45250
$ws.Cells.Item(26, 5).Value = 15
$ws.Cells.Item(26, 6).Value = 100112031
$ws.Cells.Item(26, 7).Value = "Poroto verde"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Segunda"
$ws.Cells.Item(26, 10).Value = 3500
$ws.Cells.Item(26, 11).Value = 1600
$ws.Cells.Item(26, 12).Value = 1600
$ws.Cells.Item(26, 13).Value = 1600
$ws.Cells.Item(26, 14).Value = "`$/kilo"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 1600
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = "Hortaliza"
